$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (D, K..T) that carry the per-record data for rows 3..16.
# Column letter -> 1-based column index
$cols = @{
    "D" = 4
    "K" = 11
    "L" = 12
    "M" = 13
    "N" = 14
    "O" = 15
    "P" = 16
    "Q" = 17
    "R" = 18
    "S" = 19
    "T" = 20
}

# Snapshot every value for rows 3..16 BEFORE writing anything, since the
# update re-shuffles data between rows (each row's new data comes from a
# different row's old data).
$snapshot = @{}
for ($r = 3; $r -le 16; $r++) {
    $rowData = @{}
    foreach ($colName in $cols.Keys) {
        $colIdx = $cols[$colName]
        $rowData[$colName] = $ws.Cells.Item($r, $colIdx).Value2
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (which row's original data should now live
# at the destination row).
$mapping = @{
    3  = 15
    4  = 7
    5  = 10
    6  = 11
    7  = 12
    8  = 14
    9  = 5
    10 = 6
    11 = 8
    12 = 9
    13 = 4
    14 = 13
    15 = 16
    16 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($colName in $cols.Keys) {
        $colIdx = $cols[$colName]
        $ws.Cells.Item($destRow, $colIdx).Value = $srcData[$colName]
    }
}
